# Update Sheet1 ("Employees") attendance numbers and append a new
# employee ("אייל") block, per the authoritative diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Helper: write a value as Excel *Text* (t="inlineStr"/"s"), not Number,
# while leaving the cell's style index untouched (back to default/"Normal"
# after the write so no stray numFmt/style sticks to the cell).
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# --- A2: total employee count 10 -> 11 -------------------------------
$ws.Range("A2").Value = 11

# --- Employee #3 (קרן), rows 15-19 ------------------------------------
Set-TextValue $ws.Range("B16") "6"
Set-TextValue $ws.Range("B17") "2"
Set-TextValue $ws.Range("B18") "0"

# --- Employee #4 (נדב), rows 20-24 ------------------------------------
Set-TextValue $ws.Range("B21") "6"
Set-TextValue $ws.Range("B22") "2"
Set-TextValue $ws.Range("B23") "3"

# --- Employee #5 (עידן), rows 25-29 -----------------------------------
Set-TextValue $ws.Range("B26") "6"
Set-TextValue $ws.Range("B27") "2"
Set-TextValue $ws.Range("B28") "0"

# --- Employee #6 (ליאור), rows 30-34 ----------------------------------
Set-TextValue $ws.Range("B31") "6"
Set-TextValue $ws.Range("B32") "0"
Set-TextValue $ws.Range("B33") "2"

# --- Employee #7 (נדב א), rows 35-39 ----------------------------------
Set-TextValue $ws.Range("B36") "6"
Set-TextValue $ws.Range("B37") "1"
Set-TextValue $ws.Range("B38") "0"

# --- Employee #8 (אילן), rows 40-44 -----------------------------------
Set-TextValue $ws.Range("B41") "5"
Set-TextValue $ws.Range("B42") "2"
Set-TextValue $ws.Range("B43") "0"

# --- Employee #9 (ירדן), rows 45-49 -----------------------------------
Set-TextValue $ws.Range("B46") "6"
Set-TextValue $ws.Range("B47") "0"
Set-TextValue $ws.Range("B48") "3"

# --- Employee #10 (סנד), rows 50-54 -----------------------------------
Set-TextValue $ws.Range("B51") "5"
Set-TextValue $ws.Range("B52") "4"
Set-TextValue $ws.Range("B53") "3"

# --- New employee #11 (אייל), rows 55-59 ------------------------------
$ws.Range("A55").Value = "אייל"
$ws.Range("B55").Value = 11
$ws.Range("C55").Value = "#"

$ws.Range("A56").Value = "Number of shifts"
Set-TextValue $ws.Range("B56") "4"

$ws.Range("A57").Value = "Last week Nights"
Set-TextValue $ws.Range("B57") "1"

$ws.Range("A58").Value = "Saturdays"
Set-TextValue $ws.Range("B58") "0"
$ws.Range("D58").Value = "Saturdays before"

$ws.Range("A59").Value = "Incharge"
$ws.Range("B59").Value = $true
